$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 regression-test data refresh (new candidate id / username / password,
# plus the updated numeric Candidate ID in column B).
$ws.Range("A2").Value = "test404"
$ws.Range("B2").Value = 23071131
$ws.Range("C2").Value = "narendra667"
$ws.Range("D2").Value = 'E7!$F3dy'
